# TKB_HocKy_PhongHoc.xlsx - re-run of the GA timetabling algorithm produced a
# different (still valid) room schedule for every "Phong_Tuan_*" (Room/Week)
# sheet. On each sheet the class "Tiếng Anh chuyên ngành" (was in slot
# C1 13:00-15:00 / Thứ 4) moves to slot S1 07:00-09:00 / Thứ 3, and the class
# "Kỹ năng mềm" (was in slot C1 13:00-15:00 / Thứ 5) moves to a new slot C2
# 15:00-17:00 / Thứ 6 - which requires inserting a new row 9 into the sheet.
#
# Every sheet has the identical layout, so the same sequence of operations is
# applied to each worksheet in the workbook.

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- Move the two scheduled classes to their new slots/days -----------
    # "Kỹ năng mềm" (E8, Thứ 4 / C1) -> G9 (Thứ 6 / new C2 row)
    $ws.Range("E8").Copy($ws.Range("G9"))
    # "Tiếng Anh chuyên ngành" (F8, Thứ 5 / C1) -> D8 (Thứ 3 / S1)
    $ws.Range("F8").Copy($ws.Range("D8"))

    # --- Build the rest of the new row 9, reusing row 8's formatting -------
    $ws.Range("A8").Copy($ws.Range("A9"))   # slot-label cell style (bold/border)
    $ws.Range("B8").Copy($ws.Range("B9"))   # room cell style + "R101" value

    # Blank data cells (style only, no fill) for the remaining columns of row 9
    $ws.Range("B8").Copy($ws.Range("C9"))
    $ws.Range("B8").Copy($ws.Range("D9"))
    $ws.Range("B8").Copy($ws.Range("E9"))
    $ws.Range("B8").Copy($ws.Range("F9"))
    $ws.Range("B8").Copy($ws.Range("H9"))
    $ws.Range("C9").Value = ""
    $ws.Range("D9").Value = ""
    $ws.Range("E9").Value = ""
    $ws.Range("F9").Value = ""
    $ws.Range("H9").Value = ""

    # --- Clear the two now-vacated cells in row 8 down to a blank style ----
    $ws.Range("B8").Copy($ws.Range("E8"))
    $ws.Range("E8").Value = ""
    $ws.Range("B8").Copy($ws.Range("F8"))
    $ws.Range("F8").Value = ""

    # --- Update the slot-time labels and new room label --------------------
    $ws.Range("A8").Value = "S1`n(07:00-09:00)"
    $ws.Range("A9").Value = "C2`n(15:00-17:00)"
    $ws.Range("B9").Value = "R101"

    # Row 9 should look the same as the other timetable rows
    $ws.Rows.Item(9).RowHeight = 60
}
